$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header of column A
$ws.Range("A1").Value = "h12_afstotend (mm)"

# Data rows 2-11
$A = @(22, 18, 15, 13, 12, 11, 10.5, 10, 9, 8.5)
$B = @(0.083, 0.074, 0.066, 0.057, 0.051, 0.042, 0.035, 0.027, 0.019, 0.009)
$D = @(47, 15, 11.5, 9, 8, 7, 6, 5.2, 4.7, 4)
$E = @(0.013, 0.038, 0.061, 0.09, 0.114, 0.136, 0.153, 0.187, 0.208, 0.245)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $A[$i]
    $ws.Cells.Item($row, 2).Value = $B[$i]
    $ws.Cells.Item($row, 4).Value = $D[$i]
    $ws.Cells.Item($row, 5).Value = $E[$i]
}

[void]$ws.Range("C2").Select()
